# Add a new earthquake record (row 16) to the "Sismos" sheet.
# Values taken from the target revision:
#   A16 = 16            (ID)
#   B16 = 2016-08-31 21:00:12 (Fecha) -> serial 42613.87513888889, formatted like the
#                         existing date column (built-in date/time format, numFmtId 22)
#   C16 = 4.2            (Magnitud)
#   D16 = 88             (Profundidad)
#   E16 = "Deformacion interna"          (Origen)      -> reuses existing shared string
#   F16 = "Alajuela"                     (Provincia)   -> reuses existing shared string
#   G16 = "3.6 km norta de Carrizal"     (Descripcion) -> new shared string
#   H16 = 10.1187        (Latitud)
#   I16 = -84.161         (Longitud)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sismos")

$ws.Range("A16").Value = 16
$ws.Range("B16").Value = 42613.87513888889
$ws.Range("B16").NumberFormat = "m/d/yy h:mm"
$ws.Range("C16").Value = 4.2
$ws.Range("D16").Value = 88
$ws.Range("E16").Value = "Deformacion interna"
$ws.Range("F16").Value = "Alajuela"
$ws.Range("G16").Value = "3.6 km norta de Carrizal"
$ws.Range("H16").Value = 10.1187
$ws.Range("I16").Value = -84.161
